$wb = $excel.ActiveWorkbook
$wsArr = $wb.Worksheets

# Sheet 1: ALC
$ws = $wsArr.Item(1)
$ws.Range("H2").Value = 903.0714
$ws.Range("J2").Value = 1599.4
$ws.Range("L2").Value = 1599.4
$ws.Range("N2").Value = -1825.4
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").ClearContents()
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = 0
$ws.Range("H39").Value = 826.7059
$ws.Range("I39").Value = 122
$ws.Range("J39").Value = 1833.4286
$ws.Range("K39").Value = 366
$ws.Range("L39").Value = 5500.2858
$ws.Range("M39").Value = -70
$ws.Range("N39").Value = -6092.2858
$ws.Range("H41").Value = 1157.7778
$ws.Range("I41").Value = 2312.75
$ws.Range("J41").Value = 233.8
$ws.Range("K41").Value = 2312.75
$ws.Range("L41").Value = 233.8
$ws.Range("M41").Value = -1872.75
$ws.Range("N41").Value = -1113.8
$ws.Range("H50").Value = 706.6667
$ws.Range("J50").Value = 500
$ws.Range("L50").Value = 1500
$ws.Range("N50").Value = -2450
$ws.Range("H51").Value = 7630.391
$ws.Range("J51").Value = 5700
$ws.Range("L51").Value = 5700
$ws.Range("N51").Value = -6668
$ws.Range("H53").Value = 83334344
$ws.Range("I53").Value = 333334240
$ws.Range("J53").Value = 1042.7778
$ws.Range("K53").Value = 333334240
$ws.Range("L53").Value = 1042.7778
$ws.Range("M53").Value = -333333603
$ws.Range("N53").Value = -2316.7778
$ws.Range("H58").Value = 6900
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 6900
$ws.Range("K58").Value = 0
$ws.Range("L58").ClearContents()
$ws.Range("M58").Value = 20700
$ws.Range("N58").Value = -21000
$ws.Range("H80").Value = 5604834
$ws.Range("I80").Value = 2803061.5
$ws.Range("J80").Value = 9527316
$ws.Range("K80").Value = 8409184.5
$ws.Range("L80").Value = 28581948
$ws.Range("M80").Value = -8408186.5
$ws.Range("N80").Value = -28583944
$ws.Range("H83").Value = 5604834
$ws.Range("I83").Value = 2803061.5
$ws.Range("J83").Value = 9527316
$ws.Range("K83").Value = 25227553.5
$ws.Range("L83").Value = 85745844
$ws.Range("M83").Value = -25222561.5
$ws.Range("N83").Value = -85755828
$ws.Range("H112").Value = 2604.55
$ws.Range("J112").Value = 2738.6667
$ws.Range("L112").Value = 8216.000100000001
$ws.Range("N112").Value = -10432.0001
$ws.Range("H135").Value = 1661.8788
$ws.Range("I135").Value = 504.55554
$ws.Range("K135").Value = 4540.99986
$ws.Range("M135").Value = -2005.99986

# Sheet 2: ARM
$ws = $wsArr.Item(2)
$ws.Range("H32").Value = 3457.1968
$ws.Range("I32").Value = 3469.8867
$ws.Range("J32").Value = 3373.125
$ws.Range("K32").Value = 3469.8867
$ws.Range("L32").Value = 3373.125
$ws.Range("M32").Value = -3182.8867
$ws.Range("N32").Value = -3947.125
$ws.Range("H61").Value = 14010130
$ws.Range("I61").Value = 22228572
$ws.Range("J61").Value = 1682465.6
$ws.Range("K61").Value = 22228572
$ws.Range("L61").Value = 1682465.6
$ws.Range("M61").Value = -22228360
$ws.Range("N61").Value = -1682889.6
$ws.Range("H132").Value = 3129797.5
$ws.Range("I132").Value = 4705.0835
$ws.Range("K132").Value = 14115.2505
$ws.Range("M132").Value = -11585.2505
$ws.Range("H136").Value = 14010130
$ws.Range("I136").Value = 22228572
$ws.Range("J136").Value = 1682465.6
$ws.Range("K136").Value = 66685716
$ws.Range("L136").Value = 5047396.800000001
$ws.Range("M136").Value = -66683166
$ws.Range("N136").Value = -5052496.800000001

# Sheet 3: BSM
$ws = $wsArr.Item(3)
$ws.Range("H21").Value = 34996
$ws.Range("J21").Value = 34996
$ws.Range("L21").Value = 34996
$ws.Range("N21").Value = -35468
$ws.Range("H86").Value = 357038.5
$ws.Range("I86").Value = 532338.0600000001
$ws.Range("K86").Value = 532338.0600000001
$ws.Range("M86").Value = -531215.0600000001
$ws.Range("H89").Value = 357038.5
$ws.Range("I89").Value = 532338.0600000001
$ws.Range("K89").Value = 2661690.3
$ws.Range("M89").Value = -2656074.3
$ws.Range("H133").Value = 79991
$ws.Range("J133").Value = 79991
$ws.Range("L133").Value = 79991
$ws.Range("N133").Value = -90111

# Sheet 4: CRP
$ws = $wsArr.Item(4)
$ws.Range("H31").Value = 33672650
$ws.Range("I31").Value = 55558200
$ws.Range("K31").Value = 55558200
$ws.Range("M31").Value = -55557905
$ws.Range("H34").Value = 33672650
$ws.Range("I34").Value = 55558200
$ws.Range("K34").Value = 55558200
$ws.Range("M34").Value = -55557998
$ws.Range("H52").Value = 57499.75
$ws.Range("J52").Value = 59999.668
$ws.Range("L52").Value = 59999.668
$ws.Range("N52").Value = -60587.668
$ws.Range("H86").Value = 14018.889
$ws.Range("I86").Value = 14833.75
$ws.Range("K86").Value = 14833.75
$ws.Range("M86").Value = -13710.75
$ws.Range("H89").Value = 14018.889
$ws.Range("I89").Value = 14833.75
$ws.Range("K89").Value = 74168.75
$ws.Range("M89").Value = -68552.75
$ws.Range("H93").Value = 60149
$ws.Range("J93").Value = 85223.5
$ws.Range("L93").Value = 85223.5
$ws.Range("N93").Value = -88967.5
$ws.Range("H99").Value = 45082.43
$ws.Range("I99").Value = 5665.75
$ws.Range("K99").Value = 5665.75
$ws.Range("M99").Value = -4167.75
$ws.Range("H107").Value = 5600
$ws.Range("I107").Value = 1200
$ws.Range("K107").Value = 1200
$ws.Range("M107").Value = 720
$ws.Range("H126").Value = 45082.43
$ws.Range("I126").Value = 5665.75
$ws.Range("K126").Value = 16997.25
$ws.Range("M126").Value = -14527.25
$ws.Range("H132").Value = 2840.6667
$ws.Range("I132").Value = 2810.3684
$ws.Range("J132").Value = 2955.8
$ws.Range("K132").Value = 8431.1052
$ws.Range("L132").Value = 8867.400000000001
$ws.Range("M132").Value = -5901.1052
$ws.Range("N132").Value = -13927.4
$ws.Range("H134").Value = 3305
$ws.Range("I134").Value = 3073.5
$ws.Range("J134").Value = 4346.75
$ws.Range("K134").Value = 9220.5
$ws.Range("L134").Value = 13040.25
$ws.Range("M134").Value = -6685.5
$ws.Range("N134").Value = -18110.25

# Sheet 5: CUL
$ws = $wsArr.Item(5)
$ws.Range("H3").Value = 8523
$ws.Range("I3").Value = 4388
$ws.Range("K3").Value = 13164
$ws.Range("M3").Value = -13052
$ws.Range("H97").Value = 1315.5834
$ws.Range("I97").Value = 2249
$ws.Range("J97").Value = 1128.9
$ws.Range("K97").Value = 6747
$ws.Range("L97").Value = 3386.7
$ws.Range("M97").Value = -6251
$ws.Range("N97").Value = -4378.700000000001
$ws.Range("H116").Value = 18181
$ws.Range("I116").Value = 3029
$ws.Range("J116").Value = 33333
$ws.Range("K116").Value = 9087
$ws.Range("L116").Value = 99999
$ws.Range("M116").Value = -5645
$ws.Range("N116").Value = -106883
$ws.Range("H129").Value = 4951.5293
$ws.Range("I129").Value = 2802
$ws.Range("J129").Value = 8892.333000000001
$ws.Range("K129").Value = 8406
$ws.Range("L129").Value = 26676.999
$ws.Range("M129").Value = -3406
$ws.Range("N129").Value = -36676.999
$ws.Range("H136").Value = 6575.4546
$ws.Range("J136").Value = 15650.5
$ws.Range("L136").Value = 46951.5
$ws.Range("N136").Value = -57151.5
$ws.Range("H138").Value = 9775.4375
$ws.Range("J138").Value = 11994
$ws.Range("L138").Value = 35982
$ws.Range("N138").Value = -46262

# Sheet 6: GSM
$ws = $wsArr.Item(6)
$ws.Range("H46").Value = 3500
$ws.Range("I46").Value = 3500
$ws.Range("K46").Value = 3500
$ws.Range("M46").Value = -3344
$ws.Range("H58").Value = 21405.084
$ws.Range("I58").Value = 21405.084
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 21405.084
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -21128.084

# Sheet 7: LTW
$ws = $wsArr.Item(7)
$ws.Range("H16").Value = 4079.6316
$ws.Range("I16").Value = 2209.8333
$ws.Range("J16").Value = 7285
$ws.Range("K16").Value = 2209.8333
$ws.Range("L16").Value = 7285
$ws.Range("M16").Value = -2039.8333
$ws.Range("N16").Value = -7625
$ws.Range("H46").Value = 2000
$ws.Range("I46").Value = 2000
$ws.Range("K46").Value = 2000
$ws.Range("M46").Value = -1812
$ws.Range("H132").Value = 5086.385
$ws.Range("I132").Value = 3161.5715
$ws.Range("K132").Value = 9484.7145
$ws.Range("M132").Value = -6954.7145

# Sheet 8: WVR
$ws = $wsArr.Item(8)
$ws.Range("H55").Value = 29261.75
$ws.Range("I55").Value = 9048
$ws.Range("J55").Value = 35999.668
$ws.Range("K55").Value = 9048
$ws.Range("L55").Value = 35999.668
$ws.Range("M55").Value = -8771
$ws.Range("N55").Value = -36553.668
$ws.Range("H122").Value = 4420.353
$ws.Range("I122").Value = 3998
$ws.Range("J122").Value = 4716
$ws.Range("K122").Value = 11994
$ws.Range("L122").Value = 14148
$ws.Range("M122").Value = -9544
$ws.Range("N122").Value = -19048
$ws.Range("H126").Value = 8663.808000000001
$ws.Range("I126").Value = 9102.434999999999
$ws.Range("J126").Value = 5301
$ws.Range("K126").Value = 27307.305
$ws.Range("L126").Value = 15903
$ws.Range("M126").Value = -24837.305
$ws.Range("N126").Value = -20843
$ws.Range("H131").Value = 48265.668
$ws.Range("J131").Value = 48265.668
$ws.Range("L131").Value = 48265.668
$ws.Range("N131").Value = -58345.668
$ws.Range("H132").Value = 478837.47
$ws.Range("I132").Value = 2873.2
$ws.Range("K132").Value = 8619.599999999999
$ws.Range("M132").Value = -6089.599999999999
